$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new row at 11 - this pushes the existing row 11
#    (Anthony / 44204 / 7 / Tout) down to row 12, keeping its styles.
# ------------------------------------------------------------------
$ws.Rows.Item(11).Insert()

# ------------------------------------------------------------------
# 2. Fill the new row 11 with the new entry.
# ------------------------------------------------------------------
$ws.Range("B11").Value = 44203
$ws.Range("C11").Value = "François"
$ws.Range("D11").Value = 19
$ws.Range("E11").Value = "Tout"

# ------------------------------------------------------------------
# 3. Add the brand new row 13 (merged / wrapped entry). Inserting it
#    right below row 12 makes it inherit row 12's cell styles (date
#    format on B, name style on C) the same way row 11 did above.
# ------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Range("B13").Value = 44204
$ws.Range("C13").Value = "François"

$ws.Range("D13:E13").WrapText = $true
$ws.Range("D13").Value = "27`n26"
$ws.Range("E13").Value = "Tout`nTout"
$ws.Rows.Item(13).RowHeight = 30

# Make sure the selection / active cell matches the target state.
$ws.Range("E13").Select()

# ------------------------------------------------------------------
# 4. Conditional formatting upkeep.
#    a) The catch-all rule block (C1:C10 / C17:C1048576) keeps its
#       colours but is pushed to the lowest priority.
# ------------------------------------------------------------------
$catchAll = $ws.Range("C1").FormatConditions
$catchAll.Item(1).Priority = 13
$catchAll.Item(2).Priority = 14
$catchAll.Item(3).Priority = 15
$catchAll.Item(4).Priority = 16

# b) The rule block that used to sit on C11 now belongs to the new
#    C13 row - move it there (keeps its own colours/priorities).
$oldRule = $ws.Range("C11").FormatConditions
$oldRule.Item(1).ModifyAppliesToRange($ws.Range("C13"))
$oldRule2 = $ws.Range("C13").FormatConditions
# after moving rule #1, re-fetch remaining ones still "parked" on C11
$remaining = $ws.Range("C11").FormatConditions
while ($remaining.Count -gt 0) {
    $remaining.Item(1).ModifyAppliesToRange($ws.Range("C13"))
}

# c) Re-create the per-row rule block on C11 (new row) with the same
#    four name / colour pairs used everywhere else in the sheet.
$c11 = $ws.Range("C11").FormatConditions
$r = $c11.Add(1, 3, '="Lucas"')
$r.Interior.Color = 49407
$r.Priority = 5
$r = $c11.Add(1, 3, '="Thomas"')
$r.Interior.Color = 16711935
$r.Priority = 6
$r = $c11.Add(1, 3, '="Anthony"')
$r.Interior.Color = 10498160
$r.Priority = 7
$r = $c11.Add(1, 3, '="François"')
$r.Interior.Color = 12611584
$r.Priority = 8

# d) Re-create the per-row rule block on C12 (shifted row) with the
#    same four name / colour pairs.
$c12 = $ws.Range("C12").FormatConditions
$r = $c12.Add(1, 3, '="Lucas"')
$r.Interior.Color = 49407
$r.Priority = 9
$r = $c12.Add(1, 3, '="Thomas"')
$r.Interior.Color = 16711935
$r.Priority = 10
$r = $c12.Add(1, 3, '="Anthony"')
$r.Interior.Color = 10498160
$r.Priority = 11
$r = $c12.Add(1, 3, '="François"')
$r.Interior.Color = 12611584
$r.Priority = 12
